$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.279501564648513
$ws.Range("C2").Value = 0.03613163392991225
$ws.Range("D2").Value = 0.007396171535580365
$ws.Range("E2").Value = 0.07614127728839293
$ws.Range("F2").Value = 4.461247904343708
$ws.Range("J2").Value = 0.1909334364560777
$ws.Range("K2").Value = 1.709876221727171
$ws.Range("L2").Value = 0.3211692963925259
$ws.Range("N2").Value = 4.371269726001131
$ws.Range("B3").Value = 2.243400116295362
$ws.Range("C3").Value = 0.03193831572107797
$ws.Range("D3").Value = 0.007351124521111529
$ws.Range("E3").Value = 0.07627828469593823
$ws.Range("F3").Value = 4.449932950733213
$ws.Range("J3").Value = 0.1913817980921309
$ws.Range("K3").Value = 1.67224061195239
$ws.Range("L3").Value = 0.3189368142785654
$ws.Range("N3").Value = 4.379460405894577
$ws.Range("B4").Value = 2.222451037214142
$ws.Range("C4").Value = 0.02937647600077753
$ws.Range("D4").Value = 0.007328578973993416
$ws.Range("E4").Value = 0.07638527183220489
$ws.Range("F4").Value = 4.444562865237529
$ws.Range("J4").Value = 0.1917044313802414
$ws.Range("K4").Value = 1.650109606069407
$ws.Range("L4").Value = 0.317717151666379
$ws.Range("N4").Value = 4.385348231184551
$ws.Range("B5").Value = 2.214220414130835
$ws.Range("C5").Value = 0.02833568208167492
$ws.Range("D5").Value = 0.007320688106446482
$ws.Range("E5").Value = 0.07643463369510783
$ws.Range("F5").Value = 4.44277118142368
$ws.Range("J5").Value = 0.1918478270541399
$ws.Range("K5").Value = 1.641336830198469
$ws.Range("L5").Value = 0.3172581729292219
$ws.Range("N5").Value = 4.387963406231123
$ws.Range("B6").Value = 2.212872230638027
$ws.Range("C6").Value = 0.02816304913454815
$ws.Range("D6").Value = 0.007319456466637675
$ws.Range("E6").Value = 0.07644317873540274
$ws.Range("F6").Value = 4.442497632788815
$ws.Range("J6").Value = 0.1918723582263731
$ws.Range("K6").Value = 1.639894963051887
$ws.Range("L6").Value = 0.3171842591205447
$ws.Range("N6").Value = 4.388410687540699
$ws.Range("B7").Value = 2.222338795725562
$ws.Range("C7").Value = 0.02936242670818956
$ws.Range("D7").Value = 0.00732846728985237
$ws.Range("E7").Value = 0.07638591418710483
$ws.Range("F7").Value = 4.444537095832715
$ws.Range("J7").Value = 0.191706316978653
$ws.Range("K7").Value = 1.649990298139699
$ws.Range("L7").Value = 0.317710807612805
$ws.Range("N7").Value = 4.385382626573488
$ws.Range("B8").Value = 2.266801212470284
$ws.Range("C8").Value = 0.03468306769246965
$ws.Range("D8").Value = 0.007379584404274553
$ws.Range("E8").Value = 0.07618378052857011
$ws.Range("F8").Value = 4.457019137313182
$ws.Range("J8").Value = 0.1910782164496112
$ws.Range("K8").Value = 1.696696589434879
$ws.Range("L8").Value = 0.3203682064270339
$ws.Range("N8").Value = 4.373915587436272
$ws.Range("B9").Value = 2.363650991199847
$ws.Range("C9").Value = 0.04522226193371637
$ws.Range("D9").Value = 0.007519946974237968
$ws.Range("E9").Value = 0.07596822161904448
$ws.Range("F9").Value = 4.494016150565812
$ws.Range("J9").Value = 0.1902214516060035
$ws.Range("K9").Value = 1.796049890465355
$ws.Range("L9").Value = 0.3267766505350096
$ws.Range("N9").Value = 4.35824770225183
$ws.Range("B10").Value = 2.440706717701914
$ws.Range("C10").Value = 0.05303499218543095
$ws.Range("D10").Value = 0.007646957898877815
$ws.Range("E10").Value = 0.07591935155321927
$ws.Range("F10").Value = 4.528844567540204
$ws.Range("J10").Value = 0.1898197669930433
$ws.Range("K10").Value = 1.873797700397773
$ws.Range("L10").Value = 0.33221378521948
$ws.Range("N10").Value = 4.350902866765978
$ws.Range("B11").Value = 2.477045965152172
$ws.Range("C11").Value = 0.05660550998942426
$ws.Range("D11").Value = 0.007709808262983842
$ws.Range("E11").Value = 0.07592074678237459
$ws.Range("F11").Value = 4.546353116990844
$ws.Range("J11").Value = 0.1896863337972441
$ws.Range("K11").Value = 1.910204539044912
$ws.Range("L11").Value = 0.3348454035262165
$ws.Range("N11").Value = 4.348468425253259
$ws.Range("B12").Value = 2.490991708427828
$ws.Range("C12").Value = 0.0579600206540789
$ws.Range("D12").Value = 0.007734327529195895
$ws.Range("E12").Value = 0.07592465947687366
$ws.Range("F12").Value = 4.553222735057631
$ws.Range("J12").Value = 0.189642880847213
$ws.Range("K12").Value = 1.924140485757164
$ws.Range("L12").Value = 0.3358646537956247
$ws.Range("N12").Value = 4.347677124835982
$ws.Range("B13").Value = 2.487980022870431
$ws.Range("C13").Value = 0.05766819329380724
$ws.Range("D13").Value = 0.007729015040597531
$ws.Range("E13").Value = 0.07592366649540772
$ws.Range("F13").Value = 4.551732585971166
$ws.Range("J13").Value = 0.1896519247503257
$ws.Range("K13").Value = 1.92113248157554
$ws.Range("L13").Value = 0.3356441304642317
$ws.Range("N13").Value = 4.347841735509817
$ws.Range("B14").Value = 2.478189585513917
$ws.Range("C14").Value = 0.05671689726071349
$ws.Range("D14").Value = 0.007711811121231094
$ws.Range("E14").Value = 0.07592100092923815
$ws.Range("F14").Value = 4.546913483877034
$ws.Range("J14").Value = 0.1896826171723269
$ws.Range("K14").Value = 1.911348061824498
$ws.Range("L14").Value = 0.3349288028770872
$ws.Range("N14").Value = 4.34840070663985
$ws.Range("B15").Value = 2.472216731333958
$ws.Range("C15").Value = 0.05613451976883255
$ws.Range("D15").Value = 0.007701366584825564
$ws.Range("E15").Value = 0.07591980854883218
$ws.Range("F15").Value = 4.543992839021428
$ws.Range("J15").Value = 0.1897023381909904
$ws.Range("K15").Value = 1.905374289296674
$ws.Range("L15").Value = 0.3344936008503367
$ws.Range("N15").Value = 4.348760102015348
$ws.Range("B16").Value = 2.438357735780073
$ws.Range("C16").Value = 0.05280198736571151
$ws.Range("D16").Value = 0.007642951673553
$ws.Range("E16").Value = 0.07591973467331137
$ws.Range("F16").Value = 4.527733853228057
$ws.Range("J16").Value = 0.1898294779711165
$ws.Range("K16").Value = 1.871439339904725
$ws.Range("L16").Value = 0.3320449834097872
$ws.Range("N16").Value = 4.351080224747491
$ws.Range("B17").Value = 2.417915696213527
$ws.Range("C17").Value = 0.05076185052466542
$ws.Range("D17").Value = 0.007608407873476608
$ws.Range("E17").Value = 0.07592573129110569
$ws.Range("F17").Value = 4.518185985095784
$ws.Range("J17").Value = 0.189920090254482
$ws.Range("K17").Value = 1.850887495707866
$ws.Range("L17").Value = 0.3305833350945875
$ws.Range("N17").Value = 4.352735913623121
$ws.Range("B18").Value = 2.406279045166855
$ws.Range("C18").Value = 0.04958996385917658
$ws.Range("D18").Value = 0.0075890172364943
$ws.Range("E18").Value = 0.0759314049976787
$ws.Range("F18").Value = 4.512851008487843
$ws.Range("J18").Value = 0.1899768487153715
$ws.Range("K18").Value = 1.839164406877472
$ws.Range("L18").Value = 0.3297575301814248
$ws.Range("N18").Value = 4.353773559749001
$ws.Range("B19").Value = 2.402359871959732
$ws.Range("C19").Value = 0.04919344674078729
$ws.Range("D19").Value = 0.00758253431953726
$ws.Range("E19").Value = 0.07593370849753356
$ws.Range("F19").Value = 4.51107158846321
$ws.Range("J19").Value = 0.1899968636038736
$ws.Range("K19").Value = 1.835211963855272
$ws.Range("L19").Value = 0.3294804865868031
$ws.Range("N19").Value = 4.354139540347887
$ws.Range("B20").Value = 2.420079256747556
$ws.Range("C20").Value = 0.05097886608980673
$ws.Range("D20").Value = 0.007612035718491228
$ws.Range("E20").Value = 0.07592486279741451
$ws.Range("F20").Value = 4.519186153012953
$ws.Range("J20").Value = 0.1899099642274109
$ws.Range("K20").Value = 1.853065153457067
$ws.Range("L20").Value = 0.3307373886184877
$ws.Range("N20").Value = 4.352550829361448
$ws.Range("B21").Value = 2.481060257990805
$ws.Range("C21").Value = 0.05699624941152592
$ws.Range("D21").Value = 0.007716844891795915
$ws.Range("E21").Value = 0.07592169212409594
$ws.Range("F21").Value = 4.548322469923477
$ws.Range("J21").Value = 0.1896734101464759
$ws.Range("K21").Value = 1.914217926031483
$ws.Range("L21").Value = 0.3351382958115749
$ws.Range("N21").Value = 4.348232978194076
$ws.Range("B22").Value = 2.521992329235388
$ws.Range("C22").Value = 0.06094316962297341
$ws.Range("D22").Value = 0.007789531154752893
$ws.Range("E22").Value = 0.07593934050793116
$ws.Range("F22").Value = 4.568760838496047
$ws.Range("J22").Value = 0.1895600426868107
$ws.Range("K22").Value = 1.955056042792592
$ws.Range("L22").Value = 0.3381469166724287
$ws.Range("N22").Value = 4.346172159740973
$ws.Range("B23").Value = 2.500047553940249
$ws.Range("C23").Value = 0.05883530194836339
$ws.Range("D23").Value = 0.007750357311302025
$ws.Range("E23").Value = 0.07592812109259306
$ws.Range("F23").Value = 4.557724719046007
$ws.Range("J23").Value = 0.1896167804161308
$ws.Range("K23").Value = 1.933180228689281
$ws.Range("L23").Value = 0.336529060155442
$ws.Range("N23").Value = 4.347202352062666
$ws.Range("B24").Value = 2.41910075047781
$ws.Range("C24").Value = 0.05088075019556015
$ws.Range("D24").Value = 0.007610394108612084
$ws.Range("E24").Value = 0.07592524850835858
$ws.Range("F24").Value = 4.518733496718625
$ws.Range("J24").Value = 0.189914527674528
$ws.Range("K24").Value = 1.852080346456262
$ws.Range("L24").Value = 0.3306676958023758
$ws.Range("N24").Value = 4.352634238794508
$ws.Range("B25").Value = 2.33641530885842
$ws.Range("C25").Value = 0.04235929651500214
$ws.Range("D25").Value = 0.007477739695913144
$ws.Range("E25").Value = 0.07600725391176688
$ws.Range("F25").Value = 4.482665705421596
$ws.Range("J25").Value = 0.1904131722146971
$ws.Range("K25").Value = 1.76833905723737
$ws.Range("L25").Value = 0.3249149307801673
$ws.Range("N25").Value = 4.361755317418869
